$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kenntnisse")

# Fix typo: "Power Auomate" -> "Power Automate" in A8 (A9 "Power Apps" stays as is)
$ws.Range("A8").Value = "Power Automate"

# Update the active selection to A9
$ws.Range("A9").Select()
